# Update the "PackSpecification" sheet: the `version` column (D) now holds a
# plain numeric value (1) instead of the text "1.0.0" for every data row.
$wb = $excel.ActiveWorkbook
$packSheet = $wb.Worksheets.Item("PackSpecification")

$lastRow = $packSheet.Cells.Item($packSheet.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $packSheet.Cells.Item($r, 4).Value = 1
}

# Make "PackSpecification" the active sheet/tab (it was "LetterVariant"
# before), and move its selection to F15.
$packSheet.Activate()
$packSheet.Range("F15").Select() | Out-Null
